$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Model" column before column D (Tab) ---
# This shifts old D..I (Tab, Position Field, Position, Help, Required, Copied)
# one column to the right, becoming E..J.
$ws.Columns.Item(4).Insert()

# --- Append a new "Selection Values" column after the last one ---
# Old last column I is now J; insert a fresh column at K for the new header.
$ws.Columns.Item(11).Insert()

# --- New header labels ---
# (set "Selection Values" first so it lands at sharedStrings index 9,
#  ahead of "Model" at index 10, matching the target string table order)
$ws.Range("K1").Value = "Selection Values"
$ws.Range("D1").Value = "Model"

# --- Match the "Heading 2" style used by the rest of the header row ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths matching the final layout ---
# (engine's ColumnWidth setter adds a fixed 5/6 offset vs. the raw OOXML
#  <col width> value, so subtract it to land on the exact target widths)
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666   # -> width 14
$ws.Columns.Item(11).ColumnWidth = 16.166666666666668  # -> width 17

# --- Restore the on-screen selection/active cell ---
$ws.Range("H14").Select()
